$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to
# Text format first, otherwise Excel auto-converts strings such as
# "593.29" or "13.50" into numeric values (losing trailing zeros / exact text).
$ws.Range('D2').Value = '61.172.48'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '2.929.26'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.29'
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.56'
$ws.Range('E6').Value = '  -0.73%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.02'
$ws.Range('E9').Value = '  +3.65%  '
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.441'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '33.55'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').Value = '3.415.01'
$ws.Range('E15').Value = '  +0.65%  '
$ws.Range('D16').Value = '61.108.41'
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '2.930.92'
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '432.56'
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.50'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('E21').Value = '  -0.71%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '81.90'
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.03'
$ws.Range('E24').Value = '  +1.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.20'
$ws.Range('E25').Value = '  -1.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.81'
$ws.Range('E26').Value = '  -2.61%  '
$ws.Range('E28').Value = '  -4.16%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('E30').Value = '  -1.41%  '
$ws.Range('E31').Value = '  +2.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.77'
$ws.Range('E32').Value = '  +0.32%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').Value = '0.0₃0884'
$ws.Range('E34').Value = '  +2.91%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('E37').Value = '  -4.21%  '
$ws.Range('E38').Value = '  -0.24%  '
$ws.Range('E39').Value = '  -1.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.64'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '41.97'
$ws.Range('E41').Value = '  +4.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.283'
$ws.Range('E42').Value = '  -2.46%  '
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.695.96'
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '369.22'
$ws.Range('E45').Value = '  -2.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '133.34'
$ws.Range('E46').Value = '  +2.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.83'
$ws.Range('E48').Value = '  -1.09%  '
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('E50').Value = '  -0.56%  '
$ws.Range('E51').Value = '  -0.24%  '
